# "Bonus late night session" - update the timesheet for the week of row 17
# (Sheet2) / row 18 (Sheet1): the day that used to be "9+5+8" hours became
# "9+5+10" hours, and the log entry for that day got a follow-up note about
# a late debugging session.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Extra debugging note appended to the day's log (shared string text change).
$ws2.Range("F17").Value = "building differential testing, debugging some errors (thanks to Linux for OOM killer)"

# The extra late-night hours: 8 -> 10, bumping the day total from 22 to 24.
# All of Sheet1's cumulative DONE totals (E18, F/G columns, E31, B32) and the
# chart feeding off them recompute automatically from this one formula edit.
$ws2.Range("C17").Formula = "=9+5+10"

# Leave the cursor where the author ended up after making the edit.
$ws2.Activate() | Out-Null
$ws2.Range("F17").Select() | Out-Null
